$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55; existing rows 55-57 shift down to 56-58.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new record's data.
# Columns A,B,C,E,F,G,H,I,R keep the same values as the record that used to be
# at row 55 (now at row 56), so copy them across.
$ws.Cells.Item(55, 1).Value = 5
$ws.Cells.Item(55, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(55, 3).Value = "Maule"
$ws.Cells.Item(55, 4).Value = 44509
$ws.Cells.Item(55, 5).Value = 7
$ws.Cells.Item(55, 6).Value = 100112022
$ws.Cells.Item(55, 7).Value = "Arveja Verde"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 1500
$ws.Cells.Item(55, 11).Value = 12000
$ws.Cells.Item(55, 12).Value = 12000
$ws.Cells.Item(55, 13).Value = 12000
$ws.Cells.Item(55, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(55, 15).Value = "Región del Maule"
$ws.Cells.Item(55, 16).Value = 480
$ws.Cells.Item(55, 17).Value = 25
$ws.Cells.Item(55, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same number format as the other date cells
# in column D.
$ws.Cells.Item(55, 4).NumberFormat = $ws.Cells.Item(56, 4).NumberFormat
